# Generate Report for Archive
# - Status moved from "Ready for handoff" to "In Translation" for the
#   tracked file, on the Overview roll-up sheet and on each per-locale
#   worksheet.
# - Refresh the Status column width on each sheet to fit the new
#   (shorter) status text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) and de-de (col F) status cells ---
$overview = $wb.Sheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.43
$overview.Columns.Item(6).ColumnWidth = 12.43

# --- zh-cn sheet: Status column (col C) ---
$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.43

# --- de-de sheet: Status column (col C) ---
$dede = $wb.Sheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.43
